# The template workbook opens with "Лист1" (sheet 1) already active/selected,
# matching the sheet targeted by the diff (dimension A1:M47, tabSelected="1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the worksheet's cell selection/active cell from I7 to H14.
$ws.Range("H14").Select()

# Update the page margins (stored in the xlsx as inches, but the Excel object
# model always reports/accepts them in points -> 1 inch = 72 points).
# left/right = 1.1 cm (0.43307086614173229 in), top/bottom = 1.9 cm (0.74803149606299213 in)
$ps = $ws.PageSetup
$ps.LeftMargin   = 31.181102362204726
$ps.RightMargin  = 31.181102362204726
$ps.TopMargin    = 53.85826771653544
$ps.BottomMargin = 53.85826771653544

# Set the print scaling to 90% (adds scale="90" to <pageSetup>).
$ps.Zoom = 90
